try {
    $wb = $excel.ActiveWorkbook
    $ws = $wb.ActiveSheet

    # Row 2: Resolving-Mac / Cd28 / Cd86 -> ECs (new TPM values)
    $ws.Range("A2").Value = "Resolving-Mac"
    $ws.Range("B2").Value = "Cd28"
    $ws.Range("C2").Value = "Cd86"
    $ws.Range("D2").Value = "ECs"
    $ws.Range("E2").Value = 3
    $ws.Range("F2").Value = 1
    $ws.Range("G2").Value = 5.273410666666667
    $ws.Range("H2").Value = 15.820232
    $ws.Range("I2").Value = 1
    $ws.Range("J2").Value = 1
    $ws.Range("K2").Value = 1
    $ws.Range("L2").Value = 0.3333333333333333
    $ws.Range("M2").Value = 0.02345233333333334
    $ws.Range("N2").Value = 0.070357
    $ws.Range("O2").Value = 0.0002537772683371841
    $ws.Range("P2").Value = 0.0002537772683371841
    $ws.Range("Q2").Value = 0.1236737847582222
    $ws.Range("R2").Value = 1.113064062824
    $ws.Range("S2").Value = 0.0002537772683371841
    $ws.Range("T2").Value = 0.0002537772683371841

    # Row 3: Resolving-Mac / Cd28 / Cd86 -> Resolving-Mac (new TPM values)
    $ws.Range("A3").Value = "Resolving-Mac"
    $ws.Range("B3").Value = "Cd28"
    $ws.Range("C3").Value = "Cd86"
    $ws.Range("D3").Value = "Resolving-Mac"
    $ws.Range("E3").Value = 3
    $ws.Range("F3").Value = 1
    $ws.Range("G3").Value = 5.273410666666667
    $ws.Range("H3").Value = 15.820232
    $ws.Range("I3").Value = 1
    $ws.Range("J3").Value = 1
    $ws.Range("K3").Value = 3
    $ws.Range("L3").Value = 1
    $ws.Range("M3").Value = 92.38960533333334
    $ws.Range("N3").Value = 277.168816
    $ws.Range("O3").Value = 0.9997462227316628
    $ws.Range("P3").Value = 0.9997462227316628
    $ws.Range("Q3").Value = 487.2083302539236
    $ws.Range("R3").Value = 4384.874972285312
    $ws.Range("S3").Value = 0.9997462227316628
    $ws.Range("T3").Value = 0.9997462227316628

    # Rows 4 and 5 no longer exist in the new data -- remove them entirely
    $ws.Rows("4:5").Delete()
} catch {
    Write-Output "ERROR: $_"
}
